$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "261.07"
$ws.Range("D2").ClearFormats()

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.01%"
$ws.Range("E2").ClearFormats()

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "27.11"
$ws.Range("D3").ClearFormats()

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.98%"
$ws.Range("E3").ClearFormats()

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.699"
$ws.Range("D4").ClearFormats()

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.21%"
$ws.Range("E4").ClearFormats()

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06180"
$ws.Range("D5").ClearFormats()

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "3.81%"
$ws.Range("E5").ClearFormats()

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.680"
$ws.Range("D6").ClearFormats()

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.78%"
$ws.Range("E6").ClearFormats()

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8515"
$ws.Range("D7").ClearFormats()

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.56%"
$ws.Range("E7").ClearFormats()

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9178"
$ws.Range("D8").ClearFormats()

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.01%"
$ws.Range("E8").ClearFormats()

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1411"
$ws.Range("D9").ClearFormats()

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.77%"
$ws.Range("E9").ClearFormats()

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.04686"
$ws.Range("D10").ClearFormats()

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "7.74%"
$ws.Range("E10").ClearFormats()

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07083"
$ws.Range("D11").ClearFormats()

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.07%"
$ws.Range("E11").ClearFormats()

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03155"
$ws.Range("D12").ClearFormats()

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "6.28%"
$ws.Range("E12").ClearFormats()

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09039"
$ws.Range("D13").ClearFormats()

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.80%"
$ws.Range("E13").ClearFormats()

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001534"
$ws.Range("D14").ClearFormats()

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.03%"
$ws.Range("E14").ClearFormats()

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006166"
$ws.Range("D15").ClearFormats()

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "2.11%"
$ws.Range("E15").ClearFormats()

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006085"
$ws.Range("D16").ClearFormats()

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.48%"
$ws.Range("E16").ClearFormats()

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.450"
$ws.Range("D17").ClearFormats()

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.05%"
$ws.Range("E17").ClearFormats()

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.165"
$ws.Range("D18").ClearFormats()

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.20%"
$ws.Range("E18").ClearFormats()

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.19%"
$ws.Range("E19").ClearFormats()

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.73%"
$ws.Range("E20").ClearFormats()

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.17%"
$ws.Range("E21").ClearFormats()

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.082"
$ws.Range("D22").ClearFormats()

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "6.66%"
$ws.Range("E22").ClearFormats()

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.58%"
$ws.Range("E23").ClearFormats()

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001216"
$ws.Range("D24").ClearFormats()

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.05%"
$ws.Range("E24").ClearFormats()

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.13%"
$ws.Range("E26").ClearFormats()

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001578"
$ws.Range("D27").ClearFormats()

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-7.80%"
$ws.Range("E27").ClearFormats()

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03882"
$ws.Range("D40").ClearFormats()

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.73%"
$ws.Range("E40").ClearFormats()

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1111"
$ws.Range("D41").ClearFormats()

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.26%"
$ws.Range("E41").ClearFormats()

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.004078"
$ws.Range("D42").ClearFormats()

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "7.93%"
$ws.Range("E42").ClearFormats()

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "9.46%"
$ws.Range("E43").ClearFormats()

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002184"
$ws.Range("D44").ClearFormats()

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-10.02%"
$ws.Range("E44").ClearFormats()

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005163"
$ws.Range("D45").ClearFormats()

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "0.60%"
$ws.Range("E45").ClearFormats()

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.09%"
$ws.Range("E46").ClearFormats()

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "8.06%"
$ws.Range("E47").ClearFormats()

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1635"
$ws.Range("D48").ClearFormats()

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-25.64%"
$ws.Range("E48").ClearFormats()

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002101"
$ws.Range("D49").ClearFormats()

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.09%"
$ws.Range("E49").ClearFormats()

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002001"
$ws.Range("D50").ClearFormats()

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.09%"
$ws.Range("E50").ClearFormats()

